$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 176 (shifts existing rows 176:210 down to 177:211)
$ws.Rows("176:176").Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A176").Value = 10
$ws.Range("B176").Value = "Vega Modelo de Temuco"
$ws.Range("C176").Value = "La Araucanía"
$ws.Range("D176").Value = 44188
$ws.Range("E176").Value = 9
$ws.Range("F176").Value = 100112023
$ws.Range("G176").Value = "Brócoli"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Segunda"
$ws.Range("J176").Value = 100
$ws.Range("K176").Value = 500
$ws.Range("L176").Value = 500
$ws.Range("M176").Value = 500
$ws.Range("N176").Value = "`$/unidad"
$ws.Range("O176").Value = "Provincia de Cautín"
$ws.Range("P176").Value = 500
$ws.Range("Q176").Value = 1
$ws.Range("R176").Value = "Hortaliza"
